$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://johnmoehrke.github.io/MHV-PGHD/StructureDefinition/FM-Adopted"
$wsMeta.Range("B4").Value = "Adopted"
$wsMeta.Range("B5").Value = "Adopted indication"
$wsMeta.Range("B8").Value = "2022-04-11T07:37:02-05:00"
$wsMeta.Range("B12").Value = "When this family member is Adopted."

$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("K2").Value = "Adopted indication"
$wsElem.Range("L2").Value = "When this family member is Adopted."
$wsElem.Range("Q5").Value = "https://johnmoehrke.github.io/MHV-PGHD/StructureDefinition/FM-Adopted"
